$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6587.778
$ws.Range("I62").Value = 1300
$ws.Range("J62").Value = 7248.75
$ws.Range("K62").Value = 1300
$ws.Range("L62").Value = 7248.75
$ws.Range("M62").Value = -676
$ws.Range("N62").Value = -8496.75
$ws.Range("H65").Value = 6587.778
$ws.Range("I65").Value = 1300
$ws.Range("J65").Value = 7248.75
$ws.Range("K65").Value = 6500
$ws.Range("L65").Value = 36243.75
$ws.Range("M65").Value = -3380
$ws.Range("N65").Value = -42483.75
$ws.Range("H76").Value = 5994.6113
$ws.Range("J76").Value = 7143.4287
$ws.Range("L76").Value = 7143.4287
$ws.Range("N76").Value = -7773.4287
$ws.Range("H79").Value = 5994.6113
$ws.Range("J79").Value = 7143.4287
$ws.Range("L79").Value = 7143.4287
$ws.Range("N79").Value = -9327.4287
$ws.Range("H129").Value = 23811702
$ws.Range("I129").Value = 38462930
$ws.Range("K129").Value = 115388790
$ws.Range("M129").Value = -115383790
$ws.Range("H135").Value = 997.4737
$ws.Range("I135").Value = 626.13043
$ws.Range("K135").Value = 5635.173870000001
$ws.Range("M135").Value = -3100.173870000001
$ws.Range("H137").Value = 40246.176
$ws.Range("I137").Value = 45821.875
$ws.Range("J137").Value = 3074.8333
$ws.Range("K137").Value = 137465.625
$ws.Range("L137").Value = 9224.499899999999
$ws.Range("M137").Value = -134915.625
$ws.Range("N137").Value = -14324.4999
$ws.Range("H138").Value = 2999.518
$ws.Range("I138").Value = 1941.7273
$ws.Range("J138").Value = 3381.0164
$ws.Range("K138").Value = 5825.1819
$ws.Range("L138").Value = 10143.0492
$ws.Range("M138").Value = -685.1818999999996
$ws.Range("N138").Value = -20423.0492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5388.94
$ws.Range("I32").Value = 3282.6965
$ws.Range("J32").Value = 22430.363
$ws.Range("K32").Value = 3282.6965
$ws.Range("L32").Value = 22430.363
$ws.Range("M32").Value = -2995.6965
$ws.Range("N32").Value = -23004.363
$ws.Range("H45").Value = 9528783
$ws.Range("I45").Value = 23811752
$ws.Range("J45").Value = 6804.6665
$ws.Range("K45").Value = 23811752
$ws.Range("L45").Value = 6804.6665
$ws.Range("M45").Value = -23811375
$ws.Range("N45").Value = -7558.6665
$ws.Range("H46").Value = 1885.8182
$ws.Range("J46").Value = 1638.3334
$ws.Range("L46").Value = 1638.3334
$ws.Range("N46").Value = -2276.3334
$ws.Range("H61").Value = 7983.8965
$ws.Range("I61").Value = 10071.444
$ws.Range("J61").Value = 4567.909
$ws.Range("K61").Value = 10071.444
$ws.Range("L61").Value = 4567.909
$ws.Range("M61").Value = -9859.444
$ws.Range("N61").Value = -4991.909
$ws.Range("H102").Value = 3207985
$ws.Range("I102").Value = 3790779.2
$ws.Range("K102").Value = 3790779.2
$ws.Range("M102").Value = -3789157.2
$ws.Range("I110").Value = 1544401.5
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 1544401.5
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = -1542356.5
$ws.Range("N110").Value = -7090
$ws.Range("H132").Value = 36096.9
$ws.Range("I132").Value = 10160.167
$ws.Range("J132").Value = 53388.055
$ws.Range("K132").Value = 30480.501
$ws.Range("L132").Value = 160164.165
$ws.Range("M132").Value = -27950.501
$ws.Range("N132").Value = -165224.165
$ws.Range("H136").Value = 7983.8965
$ws.Range("I136").Value = 10071.444
$ws.Range("J136").Value = 4567.909
$ws.Range("K136").Value = 30214.332
$ws.Range("L136").Value = 13703.727
$ws.Range("M136").Value = -27664.332
$ws.Range("N136").Value = -18803.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 22815480
$ws.Range("I86").Value = 43335920
$ws.Range("K86").Value = 43335920
$ws.Range("M86").Value = -43334797
$ws.Range("H89").Value = 22815480
$ws.Range("I89").Value = 43335920
$ws.Range("K89").Value = 216679600
$ws.Range("M89").Value = -216673984
$ws.Range("H134").Value = 13882.077
$ws.Range("I134").Value = 13485.117
$ws.Range("K134").Value = 40455.351
$ws.Range("M134").Value = -37920.351

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26191.826
$ws.Range("I31").Value = 11777
$ws.Range("J31").Value = 30195.945
$ws.Range("K31").Value = 11777
$ws.Range("L31").Value = 30195.945
$ws.Range("M31").Value = -11482
$ws.Range("N31").Value = -30785.945
$ws.Range("H34").Value = 26191.826
$ws.Range("I34").Value = 11777
$ws.Range("J34").Value = 30195.945
$ws.Range("K34").Value = 11777
$ws.Range("L34").Value = 30195.945
$ws.Range("M34").Value = -11575
$ws.Range("N34").Value = -30599.945
$ws.Range("H58").Value = 8697.625
$ws.Range("I58").Value = 11525.5
$ws.Range("K58").Value = 11525.5
$ws.Range("M58").Value = -11322.5
$ws.Range("H88").Value = 37031.832
$ws.Range("J88").Value = 37031.832
$ws.Range("L88").Value = 37031.832
$ws.Range("N88").Value = -37843.832
$ws.Range("H91").Value = 37031.832
$ws.Range("J91").Value = 37031.832
$ws.Range("L91").Value = 37031.832
$ws.Range("N91").Value = -39839.832
$ws.Range("H132").Value = 78287.28
$ws.Range("I132").Value = 64538.812
$ws.Range("J132").Value = 102729
$ws.Range("K132").Value = 193616.436
$ws.Range("L132").Value = 308187
$ws.Range("M132").Value = -191086.436
$ws.Range("N132").Value = -313247
$ws.Range("H133").Value = 64509.75
$ws.Range("I133").Value = 43000
$ws.Range("K133").Value = 43000
$ws.Range("M133").Value = -40470
$ws.Range("H134").Value = 7955.0605
$ws.Range("I134").Value = 5201.077
$ws.Range("K134").Value = 15603.231
$ws.Range("M134").Value = -13068.231
$ws.Range("H136").Value = 8697.625
$ws.Range("I136").Value = 11525.5
$ws.Range("K136").Value = 34576.5
$ws.Range("M136").Value = -32026.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1571.6666
$ws.Range("J122").Value = 1620
$ws.Range("L122").Value = 14580
$ws.Range("N122").Value = -19480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 49747.25
$ws.Range("J12").Value = 49747.25
$ws.Range("L12").Value = 49747.25
$ws.Range("N12").Value = -50027.25
$ws.Range("H132").Value = 9231.028
$ws.Range("I132").Value = 7713.7827
$ws.Range("J132").Value = 12139.083
$ws.Range("K132").Value = 23141.3481
$ws.Range("L132").Value = 36417.249
$ws.Range("M132").Value = -20611.3481
$ws.Range("N132").Value = -41477.249

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4362.9443
$ws.Range("I46").Value = 850.3333
$ws.Range("K46").Value = 850.3333
$ws.Range("M46").Value = -662.3333
$ws.Range("H48").Value = 29347
$ws.Range("I48").Value = 26520.5
$ws.Range("K48").Value = 26520.5
$ws.Range("M48").Value = -25859.5
$ws.Range("H100").Value = 3560.913
$ws.Range("I100").Value = 2925.25
$ws.Range("K100").Value = 2925.25
$ws.Range("M100").Value = -2384.25
$ws.Range("H122").Value = 6444.647
$ws.Range("I122").Value = 4472.375
$ws.Range("K122").Value = 13417.125
$ws.Range("M122").Value = -10967.125
$ws.Range("H132").Value = 10718.685
$ws.Range("I132").Value = 10900.541
$ws.Range("K132").Value = 32701.623
$ws.Range("M132").Value = -30171.623
$ws.Range("H136").Value = 101277.48
$ws.Range("I136").Value = 186166.36
$ws.Range("K136").Value = 558499.08
$ws.Range("M136").Value = -555949.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 100
$ws.Range("I10").Value = 100
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 100
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 69
$ws.Range("N10").ClearContents()
$ws.Range("H132").Value = 210797.9
$ws.Range("I132").Value = 7971.787
$ws.Range("J132").Value = 1572630.4
$ws.Range("K132").Value = 23915.361
$ws.Range("L132").Value = 4717891.199999999
$ws.Range("M132").Value = -21385.361
$ws.Range("N132").Value = -4722951.199999999
$ws.Range("H136").Value = 5983.722
$ws.Range("I136").Value = 6786.8213
$ws.Range("K136").Value = 20360.4639
$ws.Range("M136").Value = -17810.4639
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
